# Rename all 50 summary sheets to new postcode-based identifiers
$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ13639313",
    "summ13751513",
    "summ13869856",
    "summ14003585",
    "summ14155952",
    "summ14311264",
    "summ14462268",
    "summ14607333",
    "summ14762725",
    "summ14915811",
    "summ15081533",
    "summ15237676",
    "summ15393251",
    "summ15547890",
    "summ15696118",
    "summ15862083",
    "summ16014375",
    "summ16343386",
    "summ16521057",
    "summ16676597",
    "summ16828106",
    "summ16978563",
    "summ17129951",
    "summ17278527",
    "summ17431435",
    "summ17680499",
    "summ17832236",
    "summ17982780",
    "summ18226063",
    "summ18379691",
    "summ18532838",
    "summ18684943",
    "summ18845886",
    "summ18998884",
    "summ19173936",
    "summ19337116",
    "summ19497343",
    "summ19653345",
    "summ19801564",
    "summ19961172",
    "summ20162396",
    "summ20321977",
    "summ20460200",
    "summ20607650",
    "summ20762193",
    "summ20922147",
    "summ21063828",
    "summ21208705",
    "summ21358480",
    "summ21492925"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}

Write-Output "Renamed $($wb.Worksheets.Count) sheets"
